$wb = $excel.ActiveWorkbook

# Copy the "Swiss" sheet to the end of the workbook to create the new "Portugal" sheet
$swiss = $wb.Worksheets.Item("Swiss")
$swiss.Copy([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))

$portugal = $wb.Worksheets.Item($wb.Worksheets.Count)
$portugal.Name = "Portugal"

# Update the market-specific values on the new sheet
$portugal.Range("B2").Value = "Portugal Market"
$portugal.Range("B4").Value = "NGC-3479/T2403/T2441"

# Resize columns to fit the Portugal sheet's content (A/B/D get narrower/wider,
# C reverts to the default/standard width)
$portugal.Columns.Item(1).ColumnWidth = 24.67
$portugal.Columns.Item(2).ColumnWidth = 14.5
$portugal.Columns.Item(3).ColumnWidth = $portugal.StandardWidth
$portugal.Columns.Item(4).ColumnWidth = 16.83

# With column B narrower, the NGC code wraps onto two lines - grow those rows
$portugal.Rows.Item(3).RowHeight = 28.8
$portugal.Rows.Item(4).RowHeight = 28.8
$portugal.Rows.Item(5).RowHeight = 28.8

# Make Portugal the active sheet/tab
$portugal.Activate()
$portugal.Range("B4").Select() | Out-Null
